$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the entire content of columns B through AC between row
# 11 and row 12, while column A (the running index 9 / 10) stays
# unchanged on each row.
#
# NOTE: in this COM-interop runtime, reading a property needs to be done
# via an explicit method-call syntax, e.g. $range.Value(), otherwise a
# property-descriptor placeholder is returned instead of the real data.
# Writing, however, uses normal property assignment: $range.Value = ...

$row11Range = $ws.Range("B11:AC11")
$row12Range = $ws.Range("B12:AC12")

# Capture original values (2-D arrays) before overwriting anything.
$row11Values = $row11Range.Value()
$row12Values = $row12Range.Value()

# Swap: row 12's original data goes into row 11, and row 11's original
# data goes into row 12.
$row11Range.Value = $row12Values
$row12Range.Value = $row11Values
